# Updates productData sheet: fills in product names (column C) for the
# first several rows, and sets the productPrice / disCountPrice columns
# (D / E) for every data row to "Rs.139900" / "Rs.0" respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New product text (column C) values, row by row (row 1 is the header).
$prodText = @{
    2  = "Micromax 81cm (32) HD Ready LED TV (32T6175MHD, 2 x HDMI, 2 x USB)"
    3  = "Apple iPhone 6 (Silver, 16 GB)"
    4  = "Redmi Note 4 (Gold, 32 GB) (With 3 GB RAM)"
    5  = "Lenovo K6 Power (Silver, 32 GB)"
    6  = "Lenovo Vibe K5 Note (Gold, 32 GB)"
    7  = "null"
    8  = "null"
    9  = "null"
    10 = "null"
    11 = "null"
    12 = "null"
    13 = "null"
    14 = "null"
    15 = "null"
    16 = "null"
    17 = "null"
    18 = "null"
    19 = "null"
    20 = "null"
}
# (verified against the ground-truth diff: row 2's prodText becomes the
# Micromax TV description, rows 3-6 get their own product names, and all
# remaining rows keep the literal placeholder text "null")

for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = $prodText[$row]
    $ws.Cells.Item($row, 4).Value = "Rs.139900"
    $ws.Cells.Item($row, 5).Value = "Rs.0"
}
